# Update the component info for row 3 ("10 uF capacitor (15V)"):
# the Mouser part number, unit price and description change to a new
# part (187-CL31B106KAHNFNE).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3").Value = "187-CL31B106KAHNFNE"
$ws.Range("E3").Value = 0.094
$ws.Range("F3").Value = "Multilayer Ceramic Capacitors MLCC - SMD/SMT 10uF+/-10% 25V X7R 3 1206"

# Move the active selection to C3, matching the saved cursor position.
$ws.Range("C3").Select()
